# Apply updates to the data dictionary worksheet:
# Append three new variable rows describing composite symptom-score
# variables (upper_sx, lower_sx, systemic_sx) beneath the existing
# data, and update the selected/active cell to reflect the new
# bottom-of-sheet location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNames = @("upper_sx", "lower_sx", "systemic_sx")
$newDescs = @(
    "nose_run + nose_stuf + sneeze + throat_sr + earache (as implemented in Yan et al., 2018 PNAS https://www.pnas.org/content/115/5/1081)",
    "chest_tight + sob + cough (as implemented in Yan et al., 2018 PNAS https://www.pnas.org/content/115/5/1081)",
    "malaise + headache + mj_ache + lymph_node + sw_fever_chill (as implemented in Yan et al., 2018 PNAS https://www.pnas.org/content/115/5/1081)"
)

$startRow = 52

# Fill column A (variable names) first, then column B (descriptions),
# so new shared-string entries are interned in the same order Excel
# produced them in (names block, then descriptions block).
for ($i = 0; $i -lt $newNames.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newNames[$i]
}
for ($i = 0; $i -lt $newDescs.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $newDescs[$i]
}

# Update selection / view to mirror the resulting workbook state
$ws.Range("B55").Select()
$ws.Application.ActiveWindow.ScrollRow = 12
